$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-08 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-09 Tuesday", 2)
$d.Content.Find.Execute("61×48=2928", $true, $false, $false, $false, $false, $true, 1, $false, "96×42=4032", 2)
$d.Content.Find.Execute("89×30=2670", $true, $false, $false, $false, $false, $true, 1, $false, "87×56=4872", 2)
$d.Content.Find.Execute("49×48=2352", $true, $false, $false, $false, $false, $true, 1, $false, "61×22=1342", 2)
$d.Content.Find.Execute("82×98=8036", $true, $false, $false, $false, $false, $true, 1, $false, "28×90=2520", 2)
$d.Content.Find.Execute("31×60=1860", $true, $false, $false, $false, $false, $true, 1, $false, "35×11=385", 2)
$d.Content.Find.Execute("59×19=1121", $true, $false, $false, $false, $false, $true, 1, $false, "59×43=2537", 2)
$d.Content.Find.Execute("65×57=3705", $true, $false, $false, $false, $false, $true, 1, $false, "95×53=5035", 2)
$d.Content.Find.Execute("22×19=418", $true, $false, $false, $false, $false, $true, 1, $false, "28×49=1372", 2)
$d.Content.Find.Execute("97×41=3977", $true, $false, $false, $false, $false, $true, 1, $false, "17×92=1564", 2)
$d.Content.Find.Execute("21×86=1806", $true, $false, $false, $false, $false, $true, 1, $false, "74×61=4514", 2)
$d.Content.Find.Execute("87×68=5916", $true, $false, $false, $false, $false, $true, 1, $false, "84×32=2688", 2)
$d.Content.Find.Execute("51×85=4335", $true, $false, $false, $false, $false, $true, 1, $false, "56×17=952", 2)
$d.Content.Find.Execute("76×77=5852", $true, $false, $false, $false, $false, $true, 1, $false, "73×73=5329", 2)
$d.Content.Find.Execute("73×40=2920", $true, $false, $false, $false, $false, $true, 1, $false, "48×21=1008", 2)
$d.Content.Find.Execute("70×48=3360", $true, $false, $false, $false, $false, $true, 1, $false, "31×31=961", 2)
$d.Content.Find.Execute("39×65=2535", $true, $false, $false, $false, $false, $true, 1, $false, "75×30=2250", 2)
$d.Content.Find.Execute("89×55=4895", $true, $false, $false, $false, $false, $true, 1, $false, "35×70=2450", 2)
$d.Content.Find.Execute("56×33=1848", $true, $false, $false, $false, $false, $true, 1, $false, "27×50=1350", 2)
$d.Content.Find.Execute("81×87=7047", $true, $false, $false, $false, $false, $true, 1, $false, "65×78=5070", 2)
$d.Content.Find.Execute("41×81=3321", $true, $false, $false, $false, $false, $true, 1, $false, "46×62=2852", 2)
$d.Content.Find.Execute("21×52=1092", $true, $false, $false, $false, $false, $true, 1, $false, "80×68=5440", 2)
$d.Content.Find.Execute("73×44=3212", $true, $false, $false, $false, $false, $true, 1, $false, "52×52=2704", 2)
$d.Content.Find.Execute("66×18=1188", $true, $false, $false, $false, $false, $true, 1, $false, "68×25=1700", 2)
$d.Content.Find.Execute("65×16=1040", $true, $false, $false, $false, $false, $true, 1, $false, "13×76=988", 2)
$d.Content.Find.Execute("48×95=4560", $true, $false, $false, $false, $false, $true, 1, $false, "65×49=3185", 2)
